$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (preserving leading/trailing zeros,
# dotted-thousands style numbers, etc.) without leaving the cell tagged
# with a custom number format.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "28.516.39"
$ws.Range("E2").Value = "  +1.68%  "

Set-TextValue "D3" "1.824.78"
$ws.Range("E3").Value = "  +1.75%  "

$ws.Range("E4").Value = "  +0.18%  "

Set-TextValue "D5" "317.71"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("E6").Value = "  +0.09%  "

Set-TextValue "D7" "0.5387"
$ws.Range("E7").Value = "  +0.72%  "

$ws.Range("E8").Value = "  +6.05%  "

Set-TextValue "D9" "0.07719"
$ws.Range("E9").Value = "  +3.90%  "

$ws.Range("E10").Value = "  +2.52%  "

Set-TextValue "D11" "42.07"
$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "6.352"
$ws.Range("E12").Value = "  +3.77%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D13" "21.18"
$ws.Range("E13").Value = "  +3.02%  "

Set-TextValue "D14" "1.001"
$ws.Range("E14").Value = "  +0.08%  "

Set-TextValue "D15" "7.601"
$ws.Range("E15").Value = "  +5.02%  "

Set-TextValue "D16" "1.825.38"
$ws.Range("E16").Value = "  +2.09%  "

Set-TextValue "D17" "0.00001087"
$ws.Range("E17").Value = "  +2.80%  "

Set-TextValue "D18" "89.81"
$ws.Range("E18").Value = "  +0.79%  "

Set-TextValue "D19" "0.06575"
$ws.Range("E19").Value = "  +1.19%  "

Set-TextValue "D20" "17.75"
$ws.Range("E20").Value = "  +2.80%  "

$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("E22").Value = "  +2.92%  "

Set-TextValue "D23" "28.527.14"
$ws.Range("E23").Value = "  +1.64%  "

$ws.Range("E24").Value = "  +0.46%  "

Set-TextValue "D25" "2.258"
$ws.Range("E25").Value = "  +8.04%  "

Set-TextValue "D26" "157.92"
$ws.Range("E26").Value = "  +1.65%  "

Set-TextValue "D27" "20.77"
$ws.Range("E27").Value = "  +2.34%  "

Set-TextValue "D28" "2.453"
$ws.Range("E28").Value = "  +6.25%  "

Set-TextValue "D29" "2.037.27"
$ws.Range("E29").Value = "  +2.24%  "

Set-TextValue "D30" "124.23"
$ws.Range("E30").Value = "  +2.45%  "

Set-TextValue "D31" "1.137"
$ws.Range("E31").Value = "  +1.56%  "

Set-TextValue "D32" "0.1122"
$ws.Range("E32").Value = "  +5.42%  "

Set-TextValue "D33" "5.697"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("E34").Value = "  -0.40%  "

Set-TextValue "D35" "0.07361"
$ws.Range("E35").Value = "  +13.47%  "

Set-TextValue "D36" "0.2256"
$ws.Range("E36").Value = "  +0.37%  "

Set-TextValue "D37" "0.02351"
$ws.Range("E37").Value = "  +2.74%  "

Set-TextValue "D38" "8.962"
$ws.Range("E38").Value = "  +5.78%  "

Set-TextValue "D39" "5.200"
$ws.Range("E39").Value = "  +3.51%  "

$ws.Range("E40").Value = "  +2.17%  "

Set-TextValue "D41" "0.6289"
$ws.Range("E41").Value = "  +1.64%  "

Set-TextValue "D42" "1.191"
$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("E43").Value = "  +0.08%  "

Set-TextValue "D44" "1.401"
$ws.Range("E44").Value = "  -3.30%  "

Set-TextValue "D45" "13.47"
$ws.Range("E45").Value = "  +1.25%  "

Set-TextValue "D46" "0.5899"
$ws.Range("E46").Value = "  +1.94%  "

Set-TextValue "D47" "3.713"
$ws.Range("E47").Value = "  +1.17%  "

Set-TextValue "D48" "125.22"
$ws.Range("E48").Value = "  +0.11%  "

Set-TextValue "D49" "1.995"
$ws.Range("E49").Value = "  +3.48%  "

$ws.Range("E50").Value = "  +0.51%  "

Set-TextValue "D51" "0.06926"
$ws.Range("E51").Value = "  +1.53%  "
